$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).Style = "Normal"
}

# Numeric "Qty executed upto date" column (C) updates
$ws.Range("C8").Value = 51
$ws.Range("C9").Value = 8
$ws.Range("C10").Value = 17
$ws.Range("C11").Value = 50
$ws.Range("C12").Value = 13
$ws.Range("C13").Value = 11
$ws.Range("C14").Value = 45
$ws.Range("C15").Value = 49
$ws.Range("C16").Value = 54
$ws.Range("C17").Value = 12

# "Upto date Amount" (G) / "Amount Since prev bill" (H) text-formatted numbers
Set-TextValue "G9"  "2048.00"
Set-TextValue "G10" "8024.00"
Set-TextValue "G11" "33100.00"
Set-TextValue "G13" "1496.00"
Set-TextValue "G14" "1035.00"

Set-TextValue "G19" "45703.00"
Set-TextValue "H19" "45703.00"
Set-TextValue "G21" "45703.00"
Set-TextValue "H21" "45703.00"
